$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update C6 and C7 from "Biasa" to "Normal"
$ws1.Range("C6").Value = "Normal"
$ws1.Range("C7").Value = "Normal"

# Update selected cell on Sheet1 from E7 to C7
$ws1.Range("C7").Select()
